# Adjusted Acknowledgement on last slide
#
# The Acknowledgement paragraph on the final slide originally read as one
# run:
#   "Library Systems Developer, Fred Young, provided significant feedback
#    and contributed improvements to the look and \u201cfeel\u201d of the USB API. "
#
# It becomes four runs:
#   1) "Library Systems Developer, Fred Young who preformed peer review "
#   2) "and "
#   3) "contributed "
#   4) "improvements to the look and \u201cfeel\u201d of the USB API. "

$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)

# Find the shape/paragraph that still holds the old Acknowledgement wording.
$targetShape = $null
$targetParaIndex = -1

for ($si = 1; $si -le $s.Shapes.Count; $si++) {
    $shp = $s.Shapes.Item($si)
    if (-not $shp.HasTextFrame) { continue }
    $tr = $shp.TextFrame.TextRange
    $paraCount = $tr.Paragraphs().Count
    for ($pi = 1; $pi -le $paraCount; $pi++) {
        $candidate = $tr.Paragraphs($pi, 1)
        if ($candidate.Text -like "Library Systems Developer, Fred Young*") {
            $targetShape = $shp
            $targetParaIndex = $pi
        }
    }
}

$para = $targetShape.TextFrame.TextRange.Paragraphs($targetParaIndex, 1)

# Replace the whole paragraph's text in one assignment so it collapses back
# down to a single run, then re-split that run at the required boundaries.
$quote = [char]0x201C
$unquote = [char]0x201D
$newText = "Library Systems Developer, Fred Young who preformed peer review and contributed improvements to the look and " + $quote + "feel" + $unquote + " of the USB API. "

$whole = $para.Characters(1, $para.Text.Length)
$whole.Text = $newText

$run1Text = "Library Systems Developer, Fred Young who preformed peer review "
$run2Text = "and "
$run3Text = "contributed "
$run4Text = "improvements to the look and " + $quote + "feel" + $unquote + " of the USB API. "

$run1Start = 1
$run2Start = $run1Start + $run1Text.Length
$run3Start = $run2Start + $run2Text.Length
$run4Start = $run3Start + $run3Text.Length

# Re-assigning each span's own text forces the paragraph to carry distinct
# runs at exactly these boundaries, without changing the already-correct
# combined wording.
$para.Characters($run1Start, $run1Text.Length).Text = $run1Text
$para.Characters($run2Start, $run2Text.Length).Text = $run2Text
$para.Characters($run3Start, $run3Text.Length).Text = $run3Text
$para.Characters($run4Start, $run4Text.Length).Text = $run4Text

Write-Host "Updated paragraph:" $para.Text
